# Apply LinuxForHealth rebrand edits to the FHIR StructureDefinition workbook.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: update URL, Version, Date, Publisher ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/measure-parameter-value"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet: clear row-2 Constraint(s) and rebrand the ibm.com URLs ---
$elem = $wb.Worksheets.Item("Elements")
$elem.Range("AI2").Value = ""
$elem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/measure-parameter-value"
$elem.Range("J6").Value = "ParameterDefinition {http://linuxforhealth.org/fhir/cdm/StructureDefinition/parameter-definition-with-value}`n"
